# Scheduled-runner price/profit refresh for the Ixion Leve Profits workbook.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) for the
# affected leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Range("H8").Value = 163.4
$ws.Range("I8").Value = 163.4
$ws.Range("K8").Value = 490.2
$ws.Range("M8").Value = -351.2

# Row 32
$ws.Range("H32").Value = 1382.32
$ws.Range("J32").Value = 1508.7222
$ws.Range("L32").Value = 1508.7222
$ws.Range("N32").Value = -2160.7222

# Row 112
$ws.Range("H112").Value = 38463628
$ws.Range("I112").Value = 575
$ws.Range("J112").Value = 45456910
$ws.Range("K112").Value = 1725
$ws.Range("L112").Value = 136370730
$ws.Range("M112").Value = -617
$ws.Range("N112").Value = -136372946

# Row 116
$ws.Range("H116").Value = 5212.933
$ws.Range("I116").Value = 5395.52
$ws.Range("J116").Value = 4300
$ws.Range("K116").Value = 5395.52
$ws.Range("L116").Value = 4300
$ws.Range("M116").Value = -1953.52
$ws.Range("N116").Value = -11184

# Row 138
$ws.Range("H138").Value = 2914.4722
$ws.Range("I138").Value = 731.1579
$ws.Range("J138").Value = 5354.647
$ws.Range("K138").Value = 2193.4737
$ws.Range("L138").Value = 16063.941
$ws.Range("M138").Value = 2946.5263
$ws.Range("N138").Value = -26343.941

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1419.6316
$ws.Range("I2").Value = 871.5333000000001
$ws.Range("J2").Value = 3475
$ws.Range("K2").Value = 871.5333000000001
$ws.Range("L2").Value = 3475
$ws.Range("M2").Value = -758.5333000000001
$ws.Range("N2").Value = -3701

# Row 14
$ws.Range("H14").Value = 2000
$ws.Range("I14").Value = 2000
$ws.Range("K14").Value = 2000
$ws.Range("M14").Value = -1825

# Row 63
$ws.Range("H63").Value = 83335750
$ws.Range("I63").Value = 90911460
$ws.Range("K63").Value = 90911460
$ws.Range("M63").Value = -90910774

# Row 66
$ws.Range("H66").Value = 83335750
$ws.Range("I66").Value = 90911460
$ws.Range("K66").Value = 454557300
$ws.Range("M66").Value = -454553868

# Row 116
$ws.Range("H116").Value = 1419.6316
$ws.Range("I116").Value = 871.5333000000001
$ws.Range("J116").Value = 3475
$ws.Range("K116").Value = 871.5333000000001
$ws.Range("L116").Value = 3475
$ws.Range("M116").Value = 1422.4667
$ws.Range("N116").Value = -8063

# Row 132
$ws.Range("H132").Value = 2037.16
$ws.Range("I132").Value = 967.2807
$ws.Range("J132").Value = 5425.1113
$ws.Range("K132").Value = 2901.8421
$ws.Range("L132").Value = 16275.3339
$ws.Range("M132").Value = -371.8420999999998
$ws.Range("N132").Value = -21335.3339

# Row 137
$ws.Range("H137").Value = 50453.332
$ws.Range("J137").Value = 50453.332
$ws.Range("L137").Value = 50453.332
$ws.Range("N137").Value = -60653.332

# Row 139
$ws.Range("H139").Value = 45399.57
$ws.Range("J139").Value = 45399.57
$ws.Range("L139").Value = 45399.57
$ws.Range("N139").Value = -55679.57

# Row 140
$ws.Range("H140").Value = 54850
$ws.Range("J140").Value = 54850
$ws.Range("L140").Value = 54850
$ws.Range("N140").Value = -65210

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1419.6316
$ws.Range("I3").Value = 871.5333000000001
$ws.Range("J3").Value = 3475
$ws.Range("K3").Value = 871.5333000000001
$ws.Range("L3").Value = 3475
$ws.Range("M3").Value = -757.5333000000001
$ws.Range("N3").Value = -3703

# Row 20
$ws.Range("H20").Value = 14446.632
$ws.Range("I20").Value = 1183
$ws.Range("J20").Value = 37184.285
$ws.Range("K20").Value = 1183
$ws.Range("L20").Value = 37184.285
$ws.Range("M20").Value = -936
$ws.Range("N20").Value = -37678.285

# Row 99
$ws.Range("H99").Value = 76924344
$ws.Range("I99").Value = 100000990
$ws.Range("K99").Value = 100000990
$ws.Range("M99").Value = -99999492

# Row 107
$ws.Range("H107").Value = 1547.5
$ws.Range("J107").Value = 1150
$ws.Range("L107").Value = 1150
$ws.Range("N107").Value = -4990

$ws = $wb.Worksheets.Item("CRP")
# Row 2
$ws.Range("H2").Value = 20002.5
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 39005
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 39005
$ws.Range("M2").Value = -887
$ws.Range("N2").Value = -39231

# Row 4
$ws.Range("H4").Value = 6987.625
$ws.Range("I4").Value = 4001
$ws.Range("J4").Value = 7414.2856
$ws.Range("K4").Value = 4001
$ws.Range("L4").Value = 7414.2856
$ws.Range("M4").Value = -3889
$ws.Range("N4").Value = -7638.2856

# Row 19
$ws.Range("H19").Value = 1309.5454
$ws.Range("I19").Value = 300.625
$ws.Range("K19").Value = 300.625
$ws.Range("M19").Value = -130.625

# Row 24
$ws.Range("H24").Value = 1309.5454
$ws.Range("I24").Value = 300.625
$ws.Range("K24").Value = 300.625
$ws.Range("M24").Value = -130.625

# Row 31
$ws.Range("H31").Value = 2355.653
$ws.Range("I31").Value = 1259.5143
$ws.Range("J31").Value = 5096
$ws.Range("K31").Value = 1259.5143
$ws.Range("L31").Value = 5096
$ws.Range("M31").Value = -964.5143
$ws.Range("N31").Value = -5686

# Row 34
$ws.Range("H34").Value = 2355.653
$ws.Range("I34").Value = 1259.5143
$ws.Range("J34").Value = 5096
$ws.Range("K34").Value = 1259.5143
$ws.Range("L34").Value = 5096
$ws.Range("M34").Value = -1057.5143
$ws.Range("N34").Value = -5500

# Row 132
$ws.Range("H132").Value = 1729.174
$ws.Range("I132").Value = 1504.7561
$ws.Range("J132").Value = 3569.4
$ws.Range("K132").Value = 4514.2683
$ws.Range("L132").Value = 10708.2
$ws.Range("M132").Value = -1984.2683
$ws.Range("N132").Value = -15768.2

# Row 135
$ws.Range("H135").Value = 39800
$ws.Range("J135").Value = 39800
$ws.Range("L135").Value = 39800
$ws.Range("N135").Value = -49940

# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 89032.586
$ws.Range("I4").Value = 115631.08
$ws.Range("K4").Value = 346893.24
$ws.Range("M4").Value = -346781.24

# Row 5
$ws.Range("H5").Value = 301155.44
$ws.Range("J5").Value = 546936.8
$ws.Range("L5").Value = 1640810.4
$ws.Range("N5").Value = -1641034.4

# Row 122
$ws.Range("H122").Value = 686.1429000000001
$ws.Range("I122").Value = 344.2857
$ws.Range("J122").Value = 1028
$ws.Range("K122").Value = 3098.5713
$ws.Range("L122").Value = 9252
$ws.Range("M122").Value = -648.5713000000001
$ws.Range("N122").Value = -14152

# Row 135
$ws.Range("H135").Value = 301155.44
$ws.Range("J135").Value = 546936.8
$ws.Range("L135").Value = 4922431.2
$ws.Range("N135").Value = -4927501.2

# Row 137
$ws.Range("H137").Value = 16679176
$ws.Range("I137").Value = 11442.3
$ws.Range("J137").Value = 33346910
$ws.Range("K137").Value = 34326.89999999999
$ws.Range("L137").Value = 100040730
$ws.Range("M137").Value = -29226.89999999999
$ws.Range("N137").Value = -100050930

$ws = $wb.Worksheets.Item("GSM")
# Row 4
$ws.Range("H4").Value = 15000
$ws.Range("J4").Value = 15000
$ws.Range("L4").Value = 15000
$ws.Range("N4").Value = -15224

# Row 113
$ws.Range("H113").Value = 47620388
$ws.Range("I113").Value = 100001080
$ws.Range("J113").Value = 1575.909
$ws.Range("K113").Value = 100001080
$ws.Range("L113").Value = 1575.909
$ws.Range("M113").Value = -99998910
$ws.Range("N113").Value = -5915.909

$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 88235550
$ws.Range("J55").Value = 71428850
$ws.Range("L55").Value = 71428850
$ws.Range("N55").Value = -71429196

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 1500
$ws.Range("I122").Value = 1500
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4500
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2050
$ws.Range("N122").ClearContents()

# Row 137
$ws.Range("H137").Value = 46999.2
$ws.Range("J137").Value = 46999.2
$ws.Range("L137").Value = 46999.2
$ws.Range("N137").Value = -57199.2
